$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'trauen',
    'zielen',
    'runden',
    'greifen',
    'stellen',
    'schämen',
    'werden',
    'wechseln',
    'tauschen',
    'ärgern',
    'töten',
    'arten',
    'hören',
    'spinnen',
    'suchen',
    'stimmen',
    'scheinen',
    'machen',
    'wahren',
    'sperren',
    'lesen',
    'sterben',
    'werfen',
    'plaudern',
    'weigern',
    'heilen',
    'sinken',
    'wellen',
    'zeigen',
    'drohen',
    'klingen',
    'liefern',
    'bleiben',
    'schmecken',
    'bitten',
    'schrecken',
    'erben',
    'treffen',
    'decken',
    'heben',
    'sichern',
    'planen',
    'laufen',
    'trennen',
    'zögern',
    'atmen',
    'pfeifen',
    'boxen',
    'enden',
    'klettern',
    'spannen',
    'schwächen',
    'grüßen',
    'malen',
    'fahren',
    'testen',
    'sagen',
    'irren',
    'zünden',
    'fischen',
    'gleichen',
    'wüten',
    'tollen',
    'dringen',
    'wirken',
    'platzen',
    'bluten',
    'streifen',
    'starten',
    'reizen',
    'rufen',
    'ehren',
    'rasen',
    'siegen',
    'fordern',
    'schenken',
    'fangen',
    'dienen',
    'heulen',
    'nennen',
    'äußern',
    'achten',
    'fließen',
    'fallen',
    'stillen',
    'wehtun',
    'ändern',
    'schalten',
    'flüchten',
    'fällen',
    'regnen',
    'kaufen',
    'baden',
    'folgen',
    'tropfen',
    'schulden',
    'hindern',
    'warnen',
    'liegen',
    'kennen',
    'rechnen',
    'kehren',
    'mauern',
    'sorgen',
    'fehlen',
    'kosten',
    'wenden',
    'quälen',
    'jubeln',
    'nutzen',
    'treiben',
    'betteln',
    'kämpfen',
    'bellen',
    'klagen',
    'knarren',
    'schlafen',
    'jagen',
    'stehlen',
    'feiern',
    'schauen',
    'weichen',
    'brauen',
    'münzen',
    'gründen',
    'lohnen',
    'schicken',
    'mögen',
    'helfen',
    'wachsen',
    'warten',
    'pflanzen',
    'filmen',
    'hassen',
    'scheitern',
    'gelten',
    'binden',
    'lieben',
    'stecken',
    'foltern',
    'lehnen',
    'zählen',
    'drehen',
    'sprengen',
    'schlucken',
    'graben',
    'danken',
    'brauchen',
    'schwören',
    'fühlen',
    'loben',
    'schnellen',
    'wundern',
    'süßen',
    'seufzen',
    'tanzen',
    'stammen',
    'krachen',
    'fassen',
    'kürzen',
    'drücken',
    'räumen',
    'streichen',
    'formen',
    'schreiten',
    'spielen',
    'führen',
    'wollen',
    'wiegen',
    'freuen',
    'kichern',
    'kümmern',
    'bergen',
    'spüren',
    'stoßen',
    'lockern',
    'hauen',
    'altern',
    'wehen',
    'flehen',
    'schultern',
    'parken',
    'backen',
    'schwingen',
    'geben',
    'lügen',
    'biegen',
    'saufen',
    'rühren',
    'segnen',
    'rauchen',
    'bauen'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}
